# Remove the stray "You might also like" text that was accidentally scraped
# into the Lyrics column (column C) of several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 4, 5, 6, 7, 9, 10, 11, 12, 13, 14, 15, 16, 17, 19, 20, 22, 23, 24)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    $text = $cell.Value()
    if ($text -ne $null -and $text.Contains("You might also like")) {
        $newText = $text.Replace("You might also like", "")
        # Collapse any double space left behind by the removal, and trim a
        # trailing space if the phrase was the last thing in the cell.
        $newText = $newText -replace '  +', ' '
        $newText = $newText.TrimEnd(' ')
        $cell.Value = $newText
    }
}
